$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the source page
$ws.Name = "page_10"

# ---------------------------------------------------------------------------
# Row 1 :  " MASTER PACKAGE"  (merged A1:Z1, Calibri 13, left aligned, wrap)
# ---------------------------------------------------------------------------
$scratch = $ws.Range("AB100")
$scratch.Font.Size = 13
$scratch.Font.Bold = $false
$scratch.HorizontalAlignment = -4131   # xlLeft
$scratch.WrapText = $true
$scratch.Value = "tmp"
$scratch.Copy()

$row1 = $ws.Range("A1:Z1")
$row1.Merge() | Out-Null
$row1.PasteSpecial(-4122)              # xlPasteFormats
$ws.Range("A1").Value = " MASTER PACKAGE"
$scratch.Clear()

# ---------------------------------------------------------------------------
# Row 2 :  long header line (merged A2:Z2, Calibri 13 bold, left aligned, wrap)
# ---------------------------------------------------------------------------
$scratch = $ws.Range("AB100")
$scratch.Font.Size = 13
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4131   # xlLeft
$scratch.WrapText = $true
$scratch.Value = "tmp"
$scratch.Copy()

$row2 = $ws.Range("A2:Z2")
$row2.Merge() | Out-Null
$row2.PasteSpecial(-4122)              # xlPasteFormats
$ws.Range("A2").Value = " WesternGlove Centric8 PROD                             M12225BVS563:KONRAD                         CONSTRUCTION SKETCH DETAILS                                 MASTER"
$scratch.Clear()
